$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3   = -12.822
    14  = -12.592
    16  = -12.993
    21  = -12.842
    23  = -12.68
    25  = -11.549
    26  = -12.35
    29  = -12.378
    40  = -12.426
    53  = -12.19
    57  = -13.742
    59  = -12.679
    65  = -12.002
    69  = -10.613
    79  = -12.403
    83  = -13.345
    91  = -12.385
    93  = -10.281
    100 = -12.291
    103 = -12.155
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
